$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bulk-updated
# for every data row (rows 2-391) from 45186 (2023-09-17) to 45188 (2023-09-19).
$ws.Range("C2:C391").Value = 45188
